$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "DTaP-Hep B-IPV"
$ws.Range("A7").Value = "DTaP-Hep B-IPV"
$ws.Range("B41").Value = "Fluzone Pediatric dose Preservative-free"
$ws.Range("D43").Value = "Pack of 10 Single-dose Sprayers"
$ws.Range("D48").Value = "10 pack - 1 dose syringes No Needle"
